$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new columns at D:E, shifting the existing D:K quarterly data to F:M
$ws.Range("D1:E1").EntireColumn.Insert()

# Restore number formats for the new D:E cells (same per-row style as column F, which
# now holds what used to be column D) without creating new style entries.
for ($r = 5; $r -le 102; $r++) {
    $ws.Range("F" + $r).Copy()
    $ws.Range("D" + $r + ":E" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Populate the two new quarters of data (most recent, in columns D and E)
$ws.Range("D7").Value = 43463
$ws.Range("E7").Value = 43372
$ws.Range("D8").Value = 260300
$ws.Range("E8").Value = 248000
$ws.Range("D9").Value = 163400
$ws.Range("E9").Value = 153000
$ws.Range("D10").Value = 96900
$ws.Range("E10").Value = 95000
$ws.Range("D12").Value = "NA"
$ws.Range("E12").Value = "NA"
$ws.Range("D13").Value = 0
$ws.Range("E13").Value = 0
$ws.Range("D14").Value = 0
$ws.Range("E14").Value = 0
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("D17").Value = 215700
$ws.Range("E17").Value = 204200
$ws.Range("D18").Value = 44600
$ws.Range("E18").Value = 43800
$ws.Range("D20").Value = -300
$ws.Range("E20").Value = 0
$ws.Range("D21").Value = 53100
$ws.Range("E21").Value = 50500
$ws.Range("D22").Value = 0
$ws.Range("E22").Value = 0
$ws.Range("D23").Value = 44300
$ws.Range("E23").Value = 43800
$ws.Range("D24").Value = 9700
$ws.Range("E24").Value = 9800
$ws.Range("D25").Value = 0
$ws.Range("E25").Value = 0
$ws.Range("D26").Value = 34600
$ws.Range("E26").Value = 34000
$ws.Range("D27").Value = 34600
$ws.Range("E27").Value = 34000
$ws.Range("D28").Value = 0
$ws.Range("E28").Value = 0
$ws.Range("D29").Value = "NA"
$ws.Range("E29").Value = "NA"
$ws.Range("D30").Value = 0
$ws.Range("E30").Value = 0
$ws.Range("D31").Value = 0
$ws.Range("E31").Value = 0
$ws.Range("D32").Value = 300
$ws.Range("E32").Value = 0
$ws.Range("D33").Value = 34600
$ws.Range("E33").Value = 34000
$ws.Range("D34").Value = 0
$ws.Range("E34").Value = 0
$ws.Range("D35").Value = 34600
$ws.Range("E35").Value = 34000
$ws.Range("D38").Value = 43463
$ws.Range("E38").Value = 43372
$ws.Range("D41").Value = 43500
$ws.Range("E41").Value = 53100
$ws.Range("D42").Value = 0
$ws.Range("E42").Value = 0
$ws.Range("D43").Value = 310100
$ws.Range("E43").Value = 301000
$ws.Range("D44").Value = 270500
$ws.Range("E44").Value = 240000
$ws.Range("D45").Value = 5700
$ws.Range("E45").Value = 7800
$ws.Range("D46").Value = 629700
$ws.Range("E46").Value = 602000
$ws.Range("D47").Value = 0
$ws.Range("E47").Value = 0
$ws.Range("D48").Value = 98600
$ws.Range("E48").Value = 96800
$ws.Range("D49").Value = 97800
$ws.Range("E49").Value = 100800
$ws.Range("D50").Value = 0
$ws.Range("E50").Value = 0
$ws.Range("D51").Value = 0
$ws.Range("E51").Value = 0
$ws.Range("D52").Value = 61400
$ws.Range("E52").Value = 49500
$ws.Range("D53").Value = 0
$ws.Range("E53").Value = 0
$ws.Range("D54").Value = 887600
$ws.Range("E54").Value = 849100
$ws.Range("D57").Value = 109100
$ws.Range("E57").Value = 97900
$ws.Range("D58").Value = 0
$ws.Range("E58").Value = 0
$ws.Range("D59").Value = 32500
$ws.Range("E59").Value = 29900
$ws.Range("D60").Value = 141600
$ws.Range("E60").Value = 127900
$ws.Range("D61").Value = 0
$ws.Range("E61").Value = 0
$ws.Range("D62").Value = 18300
$ws.Range("E62").Value = 19300
$ws.Range("D63").Value = 0
$ws.Range("E63").Value = 0
$ws.Range("D64").Value = 0
$ws.Range("E64").Value = 0
$ws.Range("D65").Value = 0
$ws.Range("E65").Value = 0
$ws.Range("D66").Value = 159900
$ws.Range("E66").Value = 147100
$ws.Range("D68").Value = 0
$ws.Range("E68").Value = 0
$ws.Range("D69").Value = 0
$ws.Range("E69").Value = 0
$ws.Range("D70").Value = 0
$ws.Range("E70").Value = 0
$ws.Range("D71").Value = 0
$ws.Range("E71").Value = 0
$ws.Range("D72").Value = 679400
$ws.Range("E72").Value = 654300
$ws.Range("D73").Value = 0
$ws.Range("E73").Value = 0
$ws.Range("D74").Value = 0
$ws.Range("E74").Value = 0
$ws.Range("D75").Value = 0
$ws.Range("E75").Value = 0
$ws.Range("D76").Value = 727600
$ws.Range("E76").Value = 702000
$ws.Range("D77").Value = 0
$ws.Range("E77").Value = 0
$ws.Range("D80").Value = 43463
$ws.Range("E80").Value = 43372
$ws.Range("D81").Value = 34600
$ws.Range("E81").Value = 34000
$ws.Range("D83").Value = 8700
$ws.Range("E83").Value = 6700
$ws.Range("D84").Value = 0
$ws.Range("E84").Value = 0
$ws.Range("D85").Value = 0
$ws.Range("E85").Value = 0
$ws.Range("D86").Value = 0
$ws.Range("E86").Value = 0
$ws.Range("D87").Value = 0
$ws.Range("E87").Value = 0
$ws.Range("D88").Value = 0
$ws.Range("E88").Value = 0
$ws.Range("D89").Value = 12900
$ws.Range("E89").Value = 19600
$ws.Range("D91").Value = -8000
$ws.Range("E91").Value = -6700
$ws.Range("D92").Value = 0
$ws.Range("E92").Value = 0
$ws.Range("D93").Value = 0
$ws.Range("E93").Value = 0
$ws.Range("D94").Value = -13000
$ws.Range("E94").Value = -34200
$ws.Range("D96").Value = 0
$ws.Range("E96").Value = 0
$ws.Range("D97").Value = 0
$ws.Range("E97").Value = 0
$ws.Range("D98").Value = 0
$ws.Range("E98").Value = 0
$ws.Range("D99").Value = 0
$ws.Range("E99").Value = 0
$ws.Range("D100").Value = -9300
$ws.Range("E100").Value = -7200
$ws.Range("D101").Value = -200
$ws.Range("E101").Value = 100
$ws.Range("D102").Value = -9700
$ws.Range("E102").Value = -21700
